# Case_5_148 res_bus/vm_pu.xlsx update: "case with 380 kV done"
# Rewrites the per-bus voltage magnitude (p.u.) results for rows 2-25
# across columns B,C,D,E,F,I,J,K,L,M (columns G/H/N are unaffected).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B", "C", "D", "E", "F", "I", "J", "K", "L", "M")

# Each inner array holds the new values for columns B..M (in $cols order)
# for one worksheet row, ordered to match $rows below.
$rowData = @(
    @(1.02, 1.055467568591687, 1.058801544484407, 1.051880765516556, 1.068206579710951, 1.048956121735408, 1.060474275598555, 1.061532785797125, 1.054631040917576, 1.070912382260419),
    @(1.02, 1.056599591360931, 1.059682294013484, 1.052851651317442, 1.069228665748125, 1.049292787907898, 1.06125710648595, 1.062227470949572, 1.055414254920603, 1.071749889252791),
    @(1.02, 1.057332223478265, 1.06025229954466, 1.053480279486732, 1.06989043968057, 1.049509532388701, 1.061763214695929, 1.062676441140067, 1.055920827520811, 1.072291615656185),
    @(1.02, 1.057640255550737, 1.060491953959416, 1.053744650411985, 1.07016874884415, 1.04960038814321, 1.061975878321128, 1.062865059446148, 1.056133738110966, 1.072519310263084),
    @(1.02, 1.05769197747396, 1.060532194395399, 1.053789045057851, 1.070215484001452, 1.049615627764163, 1.062011579366741, 1.062896721737042, 1.056169483654162, 1.07255753842035),
    @(1.02, 1.057336339285416, 1.06025550172434, 1.053483811648346, 1.069894158070858, 1.049510747443641, 1.061766056727327, 1.062678961973134, 1.055923672649002, 1.072294658308146),
    @(1.02, 1.055850112817387, 1.059099176701771, 1.052208797977714, 1.068551912250932, 1.049070127603188, 1.060738927208442, 1.061767669158967, 1.054895777122226, 1.071195462320851),
    @(1.02, 1.053232220787134, 1.057062373996698, 1.049965131740276, 1.066189896155059, 1.048285267724917, 1.058925654216182, 1.060157740492031, 1.053082816477077, 1.069257036717271),
    @(1.02, 1.051487612693382, 1.055705048729484, 1.048471424522321, 1.064617371978244, 1.04775635934882, 1.057714554322025, 1.059081689918883, 1.051873048047128, 1.067963746180546),
    @(1.02, 1.050732323654786, 1.05511744252788, 1.047825124168229, 1.063936962874926, 1.047525991316725, 1.057189598627653, 1.058615091920446, 1.051348936564866, 1.067403498088884),
    @(1.02, 1.050451795459562, 1.054899198188402, 1.047585132188328, 1.063684304329437, 1.047440219868685, 1.056994525082709, 1.058441677207332, 1.051154216924433, 1.067195360109713),
    @(1.02, 1.050511968829176, 1.054946011490056, 1.047636607991911, 1.063738497021003, 1.047458627300754, 1.057036372718537, 1.058478879754021, 1.051195986819909, 1.067240008084202),
    @(1.02, 1.050709134708377, 1.055099401999578, 1.047805284872704, 1.063916076497048, 1.047518905553658, 1.05717347545434, 1.058600759435499, 1.051332841822121, 1.067386294091895),
    @(1.02, 1.050830617707934, 1.055193913443263, 1.047909221903554, 1.064025499100903, 1.047556018138039, 1.057257938150201, 1.058675840348903, 1.051417157228103, 1.06747642083702),
    @(1.02, 1.051537741563889, 1.055744048847023, 1.048514327554379, 1.064662539107706, 1.047771619723558, 1.057749382470934, 1.059112642553408, 1.051907825817623, 1.068000922834174),
    @(1.02, 1.05198133791604, 1.05608916764325, 1.048894024063793, 1.065062272426919, 1.047906500143221, 1.058057507350565, 1.05938645997335, 1.052215535845268, 1.068329863466239),
    @(1.02, 1.052240093518147, 1.056290481616669, 1.04911554138817, 1.065295478765541, 1.047985043629899, 1.058237179156596, 1.05954610927382, 1.052394991494931, 1.068521705364476),
    @(1.02, 1.052328324784786, 1.056359126516276, 1.049191081007347, 1.065375004355414, 1.048011802885298, 1.058298433749892, 1.059600534747497, 1.052456176760596, 1.068587114493644),
    @(1.02, 1.051933742850694, 1.056052138419005, 1.048853281384091, 1.065019379792397, 1.047892042195805, 1.058024453870645, 1.059357088552259, 1.052182524203586, 1.068294573692054),
    @(1.02, 1.050651073732421, 1.055054231799238, 1.047755611727881, 1.063863781678497, 1.047501160701462, 1.05713310437137, 1.058564871653206, 1.051292542576423, 1.067343217528105),
    @(1.02, 1.049844722755848, 1.054426916904261, 1.047065884109383, 1.063137649169859, 1.047254226319575, 1.056572205161164, 1.05806619804444, 1.050732736773057, 1.066744848638753),
    @(1.02, 1.050272174163859, 1.054759458088961, 1.047431481902834, 1.063522544266395, 1.047385242001423, 1.056869593316435, 1.058330608945361, 1.051029523134934, 1.06706207558698),
    @(1.02, 1.05195524897044, 1.056068870294288, 1.048871691102634, 1.065038760977162, 1.047898575521437, 1.058039389476642, 1.05937036043007, 1.052197440824929, 1.068310519694776),
    @(1.02, 1.053908890676412, 1.05758884146896, 1.050544808758573, 1.066800154908182, 1.048489171952858, 1.059394825114732, 1.060574433063002, 1.053551709064183, 1.069758343873268),
)

$rows = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowNum = $rows[$r]
    $vals = $rowData[$r]
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $ws.Range("$($cols[$c])$rowNum").Value = $vals[$c]
    }
}
